$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: title-case Spanish connector words (de, del, la, las, los, el, y)
# within a municipality / state name, e.g. "Pabellón de Arteaga" ->
# "Pabellón De Arteaga".
# ---------------------------------------------------------------------------
function TransformWords($s) {
    $words = $s.Split(" ")
    $out = @()
    foreach ($w in $words) {
        if ($w.Equals("de") -or $w.Equals("del") -or $w.Equals("la") -or $w.Equals("las") -or $w.Equals("los") -or $w.Equals("el") -or $w.Equals("y")) {
            $first = $w.Substring(0,1).ToUpper()
            $rest = $w.Substring(1)
            $out += ($first + $rest)
        } else {
            $out += $w
        }
    }
    return [string]::Join(" ", $out)
}

$lastRow = 1257

# ---------------------------------------------------------------------------
# 1. Re-title the header row (row 1) with the new machine-friendly names.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# ---------------------------------------------------------------------------
# 2. Apply the word-capitalization fix across column A (state names) and
#    column B (municipality names) for every data row, plus fix the one-off
#    "MonteMorelos" -> "Montemorelos" typo.
# ---------------------------------------------------------------------------
for ($i = 2; $i -le $lastRow; $i++) {
    $aCell = $ws.Cells.Item($i, 1)
    $aVal = $aCell.Value()
    if ($aVal -ne $null -and $aVal.GetType().Name -eq "String") {
        $newA = TransformWords($aVal)
        if (-not $newA.Equals($aVal)) {
            $aCell.Value = $newA
        }
    }

    $bCell = $ws.Cells.Item($i, 2)
    $bVal = $bCell.Value()
    if ($bVal -ne $null -and $bVal.GetType().Name -eq "String") {
        if ($bVal.Equals("MonteMorelos")) {
            $bCell.Value = "Montemorelos"
        } else {
            $newB = TransformWords($bVal)
            if (-not $newB.Equals($bVal)) {
                $bCell.Value = $newB
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Fix the 1-ULP floating point drift present on the two repeated
#    percentage values (10/10184 and 99/10184).
# ---------------------------------------------------------------------------
for ($i = 2; $i -le $lastRow; $i++) {
    $dCell = $ws.Cells.Item($i, 4)
    $dVal = $dCell.Value()
    if ($dVal -ne $null) {
        if ($dVal -eq 0.0009819324430479183) {
            $dCell.Value = 0.0009819324430479185
        } elseif ($dVal -eq 0.009721131186174391) {
            $dCell.Value = 0.009721131186174393
        }
    }
}

# ---------------------------------------------------------------------------
# 4. Drop the trailing footnote rows (1259-1263) that sit below the grand
#    total row; the sheet dimension collapses to A1:D1257 as a result.
# ---------------------------------------------------------------------------
$ws.Rows("1259:1263").Delete()
